$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be stored as text, matching original inline-string cells,
# so that values like "1.001" or "0.06515" aren't auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "30.391.99"
$ws.Range("E2").Value = "  -0.02%  "

# Row 3
$ws.Range("D3").Value = "1.875.78"
$ws.Range("E3").Value = "  -0.85%  "

# Row 4
$ws.Range("E4").Value = "  +0.08%  "

# Row 5
$ws.Range("D5").Value = "238.37"
$ws.Range("E5").Value = "  +0.01%  "

# Row 6
$ws.Range("E6").Value = "  +0.09%  "

# Row 7
$ws.Range("D7").Value = "0.4793"
$ws.Range("E7").Value = "  -0.59%  "

# Row 8
$ws.Range("E8").Value = "  -2.91%  "

# Row 9
$ws.Range("D9").Value = "0.06515"
$ws.Range("E9").Value = "  -1.41%  "

# Row 10
$ws.Range("D10").Value = "1.873.25"
$ws.Range("E10").Value = "  -0.99%  "

# Row 11
$ws.Range("D11").Value = "0.07466"

# Row 12
$ws.Range("E12").Value = "  -1.88%  "

# Row 13
$ws.Range("D13").Value = "5.078"
$ws.Range("E13").Value = "  -1.94%  "

# Row 14
$ws.Range("E14").Value = "  +0.60%  "

# Row 15
$ws.Range("D15").Value = "0.6600"
$ws.Range("E15").Value = "  -0.66%  "

# Row 16
$ws.Range("D16").Value = "30.369.15"
$ws.Range("E16").Value = "  -0.01%  "

# Row 17
$ws.Range("D17").Value = "13.29"
$ws.Range("E17").Value = "  -1.20%  "

# Row 18
$ws.Range("D18").Value = "1.001"
$ws.Range("E18").Value = "  +0.02%  "

# Row 19
$ws.Range("D19").Value = "0.000007577"
$ws.Range("E19").Value = "  -2.57%  "

# Row 20
$ws.Range("D20").Value = "2.117.10"
$ws.Range("E20").Value = "  -1.03%  "

# Row 21
$ws.Range("D21").Value = "5.301"
$ws.Range("E21").Value = "  -2.63%  "

# Row 22
$ws.Range("E22").Value = "  +0.12%  "

# Row 23
$ws.Range("D23").Value = "219.60"
$ws.Range("E23").Value = "  +13.07%  "

# Row 24
$ws.Range("E24").Value = "  +0.01%  "

# Row 25
$ws.Range("D25").Value = "9.348"
$ws.Range("E25").Value = "  -0.77%  "

# Row 26
$ws.Range("D26").Value = "167.55"
$ws.Range("E26").Value = "  +2.39%  "

# Row 27
$ws.Range("D27").Value = "18.42"
$ws.Range("E27").Value = "  +0.85%  "

# Row 28
$ws.Range("D28").Value = "1.965"
$ws.Range("E28").Value = "  +0.59%  "

# Row 29
$ws.Range("D29").Value = "1.462"
$ws.Range("E29").Value = "  +0.85%  "

# Row 30
$ws.Range("E30").Value = "  +2.26%  "

# Row 31
$ws.Range("E31").Value = "  +0.44%  "

# Row 32
$ws.Range("D32").Value = "4.027"
$ws.Range("E32").Value = "  -0.97%  "

# Row 33
$ws.Range("D33").Value = "0.05025"
$ws.Range("E33").Value = "  -1.56%  "

# Row 34
$ws.Range("E34").Value = "  +4.18%  "

# Row 35
$ws.Range("D35").Value = "0.7435"
$ws.Range("E35").Value = "  +1.38%  "

# Row 36
$ws.Range("E36").Value = "  -0.07%  "

# Row 37
$ws.Range("D37").Value = "0.01820"
$ws.Range("E37").Value = "  +0.21%  "

# Row 38
$ws.Range("E38").Value = "  -1.19%  "

# Row 39 and Row 40 (coins swapped order, plus price/volume updates)
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").Value = "2.065"
$ws.Range("E39").Value = "  -1.00%  "

$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").Value = "0.9049"
$ws.Range("E40").Value = "  -1.53%  "

# Row 41
$ws.Range("D41").Value = "106.51"
$ws.Range("E41").Value = "  -0.08%  "

# Row 42
$ws.Range("D42").Value = "5.878"
$ws.Range("E42").Value = "  -0.39%  "

# Row 43
$ws.Range("D43").Value = "0.4273"
$ws.Range("E43").Value = "  -1.48%  "

# Row 45
$ws.Range("D45").Value = "7.412"
$ws.Range("E45").Value = "  -3.44%  "

# Row 46
$ws.Range("D46").Value = "64.50"
$ws.Range("E46").Value = "  -0.68%  "

# Row 47
$ws.Range("D47").Value = "0.1273"
$ws.Range("E47").Value = "  -4.62%  "

# Row 48
$ws.Range("E48").Value = "  -5.88%  "

# Row 49
$ws.Range("D49").Value = "8.913"
$ws.Range("E49").Value = "  -1.06%  "

# Row 50
$ws.Range("D50").Value = "33.71"
$ws.Range("E50").Value = "  -1.14%  "

# Row 51
$ws.Range("D51").Value = "0.3886"
$ws.Range("E51").Value = "  +0.21%  "

